$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google")

# Create new row 6 by copying the formatting of row 5 (keeps borders/fill consistent
# with the rest of the data rows) before filling in its values.
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)

# Update existing rows 3-5: Execute column (A) Yes -> No, TestCaseID column (B) renumbered
$ws.Range("A3").Value = "No"
$ws.Range("B3").Value = "1"
$ws.Range("E3").Value = "Adidas One Grey W Gum4 Campus Core Black Footshop XgrXIq"
$ws.Range("F3").Value = "http://helper.extrapulpe.com/hcap-3-spanish.mdoc"

$ws.Range("A4").Value = "No"
$ws.Range("B4").Value = "2"

$ws.Range("A5").Value = "No"
$ws.Range("B5").Value = "3"

# Fill in the new row 6 values
$ws.Range("A6").Value = "Yes"
$ws.Range("B6").Value = "4"
$ws.Range("C6").Value = "windows"
$ws.Range("D6").Value = "pass"

$ws.Range("A6").Select()
